# Fix the "units" column on the AOCS sheet: a handful of rows had the
# wrong unit recorded next to their value (read/write functional without
# units). Correct the three mislabeled entries.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("AOCS")

$ws.Range("C3").Value = "m/s"
$ws.Range("C4").Value = "sec"
$ws.Range("C5").Value = "kg"
